$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Present" in column E (Jan 3, 2024) for each student row (rows 2-8),
# matching the existing "Present" values already set in columns C and D.
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 5).Value = "Present"
}

# Update the active selection to mirror the saved view state (F13).
$ws.Range("F13").Select()
